# Add a new "Croatia" market sheet, cloned from the existing "Turkey" sheet
# (same layout/merges/styles), placed after it as the new last + active tab,
# with the market name and ticket reference updated for Croatia.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Turkey")
$source.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2477"
